$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.367665529251099
$ws.Range("B1").Value = 2.217697858810425
$ws.Range("C1").Value = 1.633270621299744
$ws.Range("D1").Value = 1.519091963768005
$ws.Range("E1").Value = 1.49520480632782
